$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column U ("pB") corrections for control column CDOM data
$updates = @{
    2  = 0.0411267193744868
    3  = 0.137754879693119
    4  = 0.142288767825704
    5  = 0.0276048100308799
    6  = 0.0159410217066824
    7  = 0.0160903776484691
    8  = 0.0169110427241819
    9  = 0.0221480982502538
    10 = 0.0203145574392796
    11 = 0.0169569612602203
    12 = 0.0183971149351745
    13 = 0.0164253374616337
    15 = 0.0167169698854419
    16 = 0.0178324621273998
    17 = 0.12486974314811
    18 = 0.0288555822700065
    20 = 0.0165825139414517
    21 = 0.0130265704213584
    24 = 0.0145134055811091
    25 = 0.0176315629795187
    26 = 0.014983506100398
    27 = 0.0150749117426301
    28 = 0.0148432692447322
    30 = 0.138285682389224
    31 = 0.0384269901255294
    32 = 0.015682766662076
    35 = 0.0257194071524935
    36 = 0.022271862737421
    43 = 0.114751473436286
    44 = 0.049990096433465
    45 = 0.0185137536535827
    46 = 0.0165529657541441
    47 = 0.0201357834504918
    48 = 0.0274770942794866
    49 = 0.0262211304957483
    50 = 0.0198702336873087
    51 = 0.0188136883612113
    52 = 0.0210278660915613
    53 = 0.0221193188753883
    54 = 0.0205918091192174
    55 = 0.0207059825985176
    56 = 0.141527248449506
    57 = 0.0508129870452514
    58 = 0.018266266152601
    59 = 0.0190956504508666
    60 = 0.0143705557927447
    61 = 0.025486775160914
    62 = 0.0230486848737979
    63 = 0.0229091568475681
    64 = 0.0205571905932622
    65 = 0.0228505301089434
    66 = 0.0195030516029249
    67 = 0.0138982294802812
    68 = 0.0178372356582628
}

foreach ($row in $updates.Keys) {
    $ws.Range("U$row").Value = $updates[$row]
}
